$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.967.34"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "1.782.83"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.551"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -4.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.285"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0710"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "1.785.88"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.57%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "33.953.87"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.20"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.86"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.66"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.35"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0520"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("D35").Value = "1.405.27"
$ws.Range("E35").Value = "  -2.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.637"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.933"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.41"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.77%  "
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0492"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.938.50"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.03"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.81"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("E51").Value = "  -0.73%  "
